# sp_AskBrent Check ID List - add CheckID 6 (Wait Stats) and bump the
# "last updated" title string from June 23, 2013 to July 11, 2013.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10: CheckID 6, Priority 200, Wait Stats finding.
$ws.Range("A10").Value2 = 6
$ws.Range("B10").Value2 = 200
$ws.Range("C10").Value2 = "Wait Stats"
$ws.Range("D10").Value2 = "(One per wait type)"
$ws.Range("E10").Value2 = "http://BrentOzar.com/waits/(waittype)"

# Bump the document title/version string.
$ws.Range("A1").Value2 = "sp_AskBrent Check ID List - v1 July 11, 2013"

# Move the active selection back up to A2 (matches the saved view state).
$ws.Range("A2").Select()
